# Update: "Actualización desde MV -datos-"
# Adds a new quarterly data row (01-07-2021) at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row number, right after the current last data row (72).
$newRow = 73

# Write the period label as text (not as a date). Using a formula that
# evaluates to a text string and then converting it to a static value via
# PasteSpecial avoids Excel's automatic "looks like a date" conversion
# (and the style churn that comes with forcing text via NumberFormat/quote
# prefix), so the cell ends up as a plain shared-string value exactly like
# the other period cells in column A.
$labelCell = $ws.Range("A$newRow")
$labelCell.Formula = "=""01-07-2021"""
$labelCell.Copy()
$labelCell.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# Numeric indicator values for the new quarter.
$ws.Range("B$newRow").Value = 42.3
$ws.Range("C$newRow").Value = 51
$ws.Range("D$newRow").Value = 47.2
$ws.Range("E$newRow").Value = 73.7
$ws.Range("F$newRow").Value = 88
$ws.Range("G$newRow").Value = 64
$ws.Range("H$newRow").Value = 58.6
$ws.Range("I$newRow").Value = 46.9
$ws.Range("J$newRow").Value = 48.7
$ws.Range("K$newRow").Value = 50.4
